$d = $word.ActiveDocument

$d.Content.Find.Execute("LUIS GERARDO", $true, $false, $false, $false, $false,
                         $true, 1, $false, "yrty", 2)

$d.Content.Find.Execute("luisgerardocazares@gmail.com", $true, $false, $false, $false, $false,
                         $true, 1, $false, "yrtyr", 2)

$d.Content.Find.Execute("5858", $true, $false, $false, $false, $false,
                         $true, 1, $false, "12112", 2)

$d.Content.Find.Execute("5525003847", $true, $false, $false, $false, $false,
                         $true, 1, $false, "yrty", 2)
